$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row 3 for the "AminoAcid" structure. This shifts the
# existing rows 3-13 down to 4-14 (values, formulas and formatting all
# move with the insert). Range objects captured in PowerShell are NOT
# "live" references that track the shift, so from here on we always
# address cells using their post-insert (final) row numbers.
# ------------------------------------------------------------------
$ws.Rows("3:3").Insert()

# "Format donor" cells, addressed where their content now lives after
# the insert above: B4 still carries the plain "Oui"/"S/O" look (the
# "Bueno" style) and C10 still carries the red "TODO" look (the
# "Incorrecto" style). Reusing these via copy/paste-special keeps the
# workbook's existing style table intact instead of Excel fabricating
# new (duplicate) style records.
$buenoDonor = $ws.Range("B4")
$todoDonor  = $ws.Range("C10")

# Stamp the whole new row with the standard "Bueno" look first ...
$buenoDonor.Copy() | Out-Null
$ws.Range("B3:L3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ... then overlay the red "TODO" style on the two cells that need it.
$todoDonor.Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Fill in the values for the new "AminoAcid" row.
# ------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "AminoAcid"
$ws.Cells.Item(3, 2).Value = "Oui"
$ws.Cells.Item(3, 3).Value = "S/O"
$ws.Cells.Item(3, 4).Value = "S/O"
$ws.Cells.Item(3, 5).Value = "S/O"
$ws.Cells.Item(3, 6).Value = "TODO"
$ws.Cells.Item(3, 7).Value = "S/O"
$ws.Cells.Item(3, 8).Value = "S/O"
$ws.Cells.Item(3, 9).Value = "S/O"
$ws.Cells.Item(3, 10).Value = "S/O"
$ws.Cells.Item(3, 11).Value = "S/O"
$ws.Cells.Item(3, 12).Value = "TODO"

# ------------------------------------------------------------------
# Update the interface coverage that changed for the structures that
# were already in the sheet (rows shifted down by one: PointD is now
# row 10, PointM row 11, PolarVector row 12, TwoDVector row 13).
# ------------------------------------------------------------------

# PointD (row 10): IEquatable<T> and IComparable<T> generic coverage
# went from TODO to Oui -> switch cell style back to "Bueno" too.
$buenoDonor.Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(10, 3).Value = "Oui"
$ws.Cells.Item(10, 5).Value = "Oui"

# PointM (row 11): same fix as PointD.
$buenoDonor.Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(11, 3).Value = "Oui"
$ws.Cells.Item(11, 5).Value = "Oui"

# PolarVector (row 12): IComparable<T> non-generic coverage went from
# S/O to Oui (style already "Bueno", only the text changes).
$ws.Cells.Item(12, 4).Value = "Oui"

# TwoDVector (row 13): IEquatable<T> coverage went from TODO to S/O
# (style back to "Bueno"), and IComparable<T> non-generic coverage
# went from S/O to Oui.
$buenoDonor.Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(13, 3).Value = "S/O"
$ws.Cells.Item(13, 4).Value = "Oui"

# ------------------------------------------------------------------
# Match the author's last selection before saving.
# ------------------------------------------------------------------
$ws.Range("D15").Select()
